# Update "想去人数" (want-to-go count) values on two sheets to reflect
# newly generated output data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 282
$ws1.Range("F4").Value = 950
$ws1.Range("F6").Value = 53

# Sheet "全部类型" (All Types)
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F4").Value = 282
$ws2.Range("F5").Value = 950
$ws2.Range("F7").Value = 53
